$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "Part 1" -> "Part 2"
$ws.Name = "Part 2"

# --- Inputs block (rows 1-5) ---
# Row 1 (Inputs / Values headers) is unchanged.

# Column A labels for the inputs block
$ws.Range("A2").Value = "Arrival Rate"
$ws.Range("A3").Value = "Service Rate"
$ws.Range("A4").Value = "P(W > 0) Less Than"
$ws.Range("A5").Value = "E(W) Less Than"

# Column B values for the inputs block (stored as text, matching the
# workbook's existing convention of keeping numeric-looking values as text)
$ws.Range("B2").Value = "'100.0"
$ws.Range("B3").Value = "'101.0"
$ws.Range("B4").Value = "'0.2"
$ws.Range("B5").Value = "'1.0"

# --- Results block (rows 6-9) ---
# Row 6 "Results" header is unchanged.

# Number of Servers stays, its value changes
$ws.Range("A7").Value = "Number of Servers"
$ws.Range("B7").Value = "'3"

# Replace "Actual P(block)" with "E(S)"
$ws.Range("A8").Value = "E(S)"
$ws.Range("B8").Value = "'0.0033003300330033004"

# New row: E(N)
$ws.Range("A9").Value = "E(N)"
$ws.Range("B9").Value = "'0.49261083743842365"
